# Add a new "Contrib_name" column to the metadata table (between Contrib_ORCID
# and Pub_ROR), populate it with the contributor's name, write the fresh DOI
# placeholder into the title/doi row, and stretch the table down by two
# (currently empty, pre-formatted) rows ready for the minted DOI write-back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room for the new "Contrib_name" column: insert a blank column at E,
#    which pushes Pub_ROR/url/publication_year/doi/title (old E:I) to F:J.
$ws.Columns("E:E").Insert() | Out-Null

# 2) Populate the new column.
$ws.Range("E1").Value = "Contrib_name"
$ws.Range("E2").Value = "Rossow, Nick"

# 3) The ORCID isn't sufficient on its own for the published DOI - pair it
#    with the contributor's name using the same fill/border styling as the
#    other new-contributor cell (Contrib_ORCID, D2).
$ws.Range("E2").Style = $ws.Range("D2").Style

# 4) Column widths: give the newly-visible/previously-default columns
#    explicit widths matching the rest of the sheet's layout.
$ws.Columns("A:A").ColumnWidth = 14.83203125
$ws.Columns("B:B").ColumnWidth = 19.83203125
$ws.Columns("C:C").ColumnWidth = 15
$ws.Columns("D:D").ColumnWidth = 17.6640625
$ws.Columns("E:E").ColumnWidth = 17.6640625

# 5) Resize the table to include the new column (F1:J2 after the shift).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("F1:J2")) | Out-Null

# 6) Re-point the hyperlink at the url cell's new location (G2).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G2"), "https://phenomicsaustralia.org.au/project") | Out-Null

# 7) New version 3 draft - write back the minted DOI placeholder.
$ws.Range("J2").Value = "fresh123"

# 8) Leave two blank, pre-formatted rows under the data row so the DOI
#    write-back script has somewhere to land future rows.
$ws.Range("F2:J2").Copy() | Out-Null
$ws.Range("F3:J4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 9) Reflect the current working selection/scroll position.
$ws.Range("F7").Select() | Out-Null
